$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("O2").Value = 1.2
$ws.Range("P2").Value = 4.5
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 2.2

# Row 4
$ws.Range("G4").Value = 1.83
$ws.Range("H4").Value = 2.9
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 2.63
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5
$ws.Range("Z4").Value = 15
$ws.Range("AB4").Value = 51
$ws.Range("AH4").Value = 9.5
$ws.Range("AI4").Value = 26
$ws.Range("AJ4").Value = 21
$ws.Range("AK4").Value = 67
$ws.Range("AL4").Value = 51
$ws.Range("AN4").Value = 3.5
$ws.Range("AO4").Value = 11
$ws.Range("AX4").Value = 6.5
$ws.Range("AY4").Value = 34

# Row 6
$ws.Range("G6").Value = 1.9
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 5
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 6
$ws.Range("X6").Value = 7
$ws.Range("Z6").Value = 15
$ws.Range("AE6").Value = 23
$ws.Range("AH6").Value = 8.5
$ws.Range("AI6").Value = 23
$ws.Range("AJ6").Value = 19
$ws.Range("AO6").Value = 11
$ws.Range("AP6").Value = 29
$ws.Range("AX6").Value = 6.5
$ws.Range("AY6").Value = 34

# Row 7
$ws.Range("G7").Value = 2.8
$ws.Range("I7").Value = 2.7
$ws.Range("J7").Value = 3.5
$ws.Range("L7").Value = 3.5
$ws.Range("Q7").Value = 2.4
$ws.Range("R7").Value = 1.53
$ws.Range("X7").Value = 12
$ws.Range("Z7").Value = 29
$ws.Range("AC7").Value = 7
$ws.Range("AF7").Value = 51
$ws.Range("AG7").Value = 401
$ws.Range("AH7").Value = 7.5
$ws.Range("AI7").Value = 12
$ws.Range("AJ7").Value = 11
$ws.Range("AK7").Value = 29
$ws.Range("AL7").Value = 26
$ws.Range("AM7").Value = 41
$ws.Range("AN7").Value = 4.5
$ws.Range("AO7").Value = 17
$ws.Range("AR7").Value = 81
$ws.Range("AW7").Value = 126
$ws.Range("AX7").Value = 4.5
$ws.Range("AY7").Value = 17

# Row 8
$ws.Range("G8").Value = 2.35
$ws.Range("H8").Value = 2.55
$ws.Range("J8").Value = 3.5
$ws.Range("AI8").Value = 17

# Row 9
$ws.Range("L9").Value = 6
$ws.Range("Q9").Value = 2.5
$ws.Range("R9").Value = 1.5
$ws.Range("AQ9").Value = 34
$ws.Range("AY9").Value = 34
$ws.Range("BB9").Value = 201

# Row 12
$ws.Range("G12").Value = 1.57
$ws.Range("I12").Value = 5.75
$ws.Range("L12").Value = 5.5
$ws.Range("Q12").Value = 1.73
$ws.Range("R12").Value = 2.08
$ws.Range("S12").Value = 1.33
$ws.Range("T12").Value = 3.25
$ws.Range("U12").Value = 1.8
$ws.Range("V12").Value = 1.91
$ws.Range("W12").Value = 7.5
$ws.Range("Z12").Value = 12
$ws.Range("AB12").Value = 23
$ws.Range("AH12").Value = 15
$ws.Range("AJ12").Value = 17
$ws.Range("AK12").Value = 51
$ws.Range("AN12").Value = 3.6
$ws.Range("AT12").Value = 3.25
$ws.Range("BC12").Value = 201

# Row 15
$ws.Range("Q15").Value = 1.57
$ws.Range("R15").Value = 2.35

# Row 28
$ws.Range("J28").Value = 3.5
$ws.Range("L28").Value = 3.4
$ws.Range("Z28").Value = 29
$ws.Range("AK28").Value = 26
$ws.Range("AS28").Value = 201
$ws.Range("AY28").Value = 15
$ws.Range("BC28").Value = 201

